# "Updated project plan (issue 18) and corrected spelling mistakes in group list"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct spelling mistakes in the group list:
#   "Josh Doyle "  -> "Joshua Doyle "
#   "Olver Earl"   -> "Oliver Earl"
$ws.Range("A5").Value2 = "Joshua Doyle "
$ws.Range("A8").Value2 = "Oliver Earl"

# Leave the active selection on C4, as last left in the updated project plan.
$ws.Range("C4").Select()
